$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84, shifting existing rows 84:166 down to 85:167
$ws.Rows("84").Insert()

# Populate the newly inserted row 84 with the new data record
$ws.Range("A84").Value = 6
$ws.Range("B84").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C84").Value = "Metropolitana"
$ws.Range("D84").Value = 44512
$ws.Range("E84").Value = 13
$ws.Range("F84").Value = 100112022
$ws.Range("G84").Value = "Arveja Verde"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 200
$ws.Range("K84").Value = 10000
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = 11200
$ws.Range("N84").Value = '$/saco 25 kilos'
$ws.Range("O84").Value = "Región Metropolitana"
$ws.Range("P84").Value = 448
$ws.Range("Q84").Value = 25
$ws.Range("R84").Value = "Hortaliza"
